$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("839:840").Insert()

$ws.Range("A839").Value = 9
$ws.Range("B839").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C839").Value = "Metropolitana"
$ws.Range("D839").Value = 45013
$ws.Range("E839").Value = 13
$ws.Range("F839").Value = 100112040
$ws.Range("G839").Value = "Cilantro"
$ws.Range("H839").Value = "Sin especificar"
$ws.Range("I839").Value = "Primera"
$ws.Range("J839").Value = 70
$ws.Range("K839").Value = 7000
$ws.Range("L839").Value = 7000
$ws.Range("M839").Value = 7000
$ws.Range("N839").Value = "$/caja 36 atados"
$ws.Range("O839").Value = "Región Metropolitana"
$ws.Range("P839").Value = 194
$ws.Range("Q839").Value = 36
$ws.Range("R839").Value = "Hortaliza"

$ws.Range("A840").Value = 9
$ws.Range("B840").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C840").Value = "Metropolitana"
$ws.Range("D840").Value = 45013
$ws.Range("E840").Value = 13
$ws.Range("F840").Value = 100112040
$ws.Range("G840").Value = "Cilantro"
$ws.Range("H840").Value = "Sin especificar"
$ws.Range("I840").Value = "Primera"
$ws.Range("J840").Value = 160
$ws.Range("K840").Value = 12000
$ws.Range("L840").Value = 13000
$ws.Range("M840").Value = 12500
$ws.Range("N840").Value = "$/docena de atados"
$ws.Range("O840").Value = "Región Metropolitana"
$ws.Range("P840").Value = 4167
$ws.Range("Q840").Value = 3
$ws.Range("R840").Value = "Hortaliza"
